# "Add files via upload" - adds a new worksheet "Tabelle1" right after the
# first sheet ("AZ MIF 726 "), scribbles a few test values on it, and adds a
# stray "Test123" label underneath the "AZ1518 Resist" sheet's table.

$wb = $excel.ActiveWorkbook

# The sheet that currently sits first ("AZ MIF 726 ") is used as the anchor
# for where the new sheet gets inserted.
$firstSheet = $wb.Worksheets.Item(1)

# The "AZ1518 Resist" sheet picks up a new, disconnected label further down
# the sheet (row 21). Write this first so the shared-string table gets the
# "Test123" entry before the strings used on the new sheet.
$az1518 = $wb.Worksheets.Item("AZ1518 Resist")
$az1518.Range("B21").Value = "Test123"

# Insert the new sheet right after "AZ MIF 726 " and rename it.
$newSheet = $wb.Worksheets.Add($null, $firstSheet)
$newSheet.Name = "Tabelle1"

# Match this workbook's usual (2 cm) top/bottom page margins, like every
# other sheet in the file already uses.
$newSheet.PageSetup.TopMargin = 56.692913399999995
$newSheet.PageSetup.BottomMargin = 56.692913399999995

# Scattered test cells on the new sheet (order chosen to reproduce the
# shared-string insertion order: e, we, w).
$newSheet.Range("B14").Value = "e"
$newSheet.Range("D28").Value = "we"
$newSheet.Range("G12").Value = "w"
$newSheet.Range("D16").Value = "w"
$newSheet.Range("G21").Value = "w"
$newSheet.Range("C25").Value = "e"

# Restore/update the selections on each sheet, leaving the new "Tabelle1"
# sheet as the active (selected) tab.
$firstSheet.Range("D32").Select()
$newSheet.Range("F19").Select()
